$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '26.271.91'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.679.08'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.42'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5273'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.90%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2705'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06486'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07540'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('D12').Value = '1.690.18'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.528'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5808'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008508'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.62'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '26.300.71'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.926'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.43'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.196'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.809'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1244'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.80'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06541'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.599'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.589'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.032'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6232'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.402'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.730'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.452'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.75%  '
$ws.Range('D39').Value = '1.112.17'
$ws.Range('E39').Value = '  +2.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01625'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8763'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.014'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.74'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '1.830.33'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.195'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.007'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05272'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.088'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.21%  '
$ws.Range('E51').Value = '  -0.03%  '
